# DDAS_Upload_Template.xlsx: rework the single upload-template sheet so the
# header row captures PI + repeated Sub-Investigator blocks instead of a
# single Investigator/Role block (per "moving files from pat's repo to admin
# repo").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header text, columns A1:AC1 (A:H unique, I:AC = 7x repeat of the
#    Sub Investigator / Sub Investigator ML# / SI Qualification triple).
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "PI Name"
$ws.Range("B1").Value = "PI Medical license #"
$ws.Range("C1").Value = "PI Qualification"
$ws.Range("D1").Value = "Project Number"
$ws.Range("E1").Value = "Sponsor Protocol #"
$ws.Range("F1").Value = "Institute Name"
$ws.Range("G1").Value = "Address"
$ws.Range("H1").Value = "Country"
$ws.Range("I1").Value = "Sub Investigator"
$ws.Range("J1").Value = "Sub Investigator ML#"
$ws.Range("K1").Value = "SI Qualification"
$ws.Range("L1").Value = "Sub Investigator"
$ws.Range("M1").Value = "Sub Investigator ML#"
$ws.Range("N1").Value = "SI Qualification"
$ws.Range("O1").Value = "Sub Investigator"
$ws.Range("P1").Value = "Sub Investigator ML#"
$ws.Range("Q1").Value = "SI Qualification"
$ws.Range("R1").Value = "Sub Investigator"
$ws.Range("S1").Value = "Sub Investigator ML#"
$ws.Range("T1").Value = "SI Qualification"
$ws.Range("U1").Value = "Sub Investigator"
$ws.Range("V1").Value = "Sub Investigator ML#"
$ws.Range("W1").Value = "SI Qualification"
$ws.Range("X1").Value = "Sub Investigator"
$ws.Range("Y1").Value = "Sub Investigator ML#"
$ws.Range("Z1").Value = "SI Qualification"
$ws.Range("AA1").Value = "Sub Investigator"
$ws.Range("AB1").Value = "Sub Investigator ML#"
$ws.Range("AC1").Value = "SI Qualification"

# ---------------------------------------------------------------------
# 2. Drop the old trailing blank header/body columns (AD:BF on row 1, AD on
#    row 2) so the used range shrinks back down to A1:AC2.
# ---------------------------------------------------------------------
$ws.Range("AD1:BF1").Clear()
$ws.Range("AD2").Clear()

# ---------------------------------------------------------------------
# 3. Re-stripe formatting. Every header cell (row 1) carries the bold
#    centered "s=1" style; row 2 (the blank data row under the headers)
#    keeps a handful of distinct alignment styles per column group.
#
#    A2 (s=6), E2 (s=5), I2 (s=2) and T2 (s=4) never change value or style
#    in this pass, so they're used as the copy-source anchors below and are
#    read before anything that *does* move is touched (G2 needs H2's
#    original "s=3" before H2 itself flips to "s=2").
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A1:AC1").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("B2").PasteSpecial(-4122)

$ws.Range("I2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$ws.Range("I2").Copy()
$ws.Range("F2").PasteSpecial(-4122)

$ws.Range("I2").Copy()
$ws.Range("H2").PasteSpecial(-4122)

$ws.Range("T2").Copy()
$ws.Range("S2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Column widths: shift to line up with the new layout (best effort —
#    the COM width setter here snaps to whole-pixel increments).
# ---------------------------------------------------------------------
$widths = @(18.5703125,18.5703125,23.85546875,14.85546875,23.85546875,14.42578125,32.28515625,16.85546875,26,19.5703125,17.85546875,15.42578125,19.5703125,19.28515625,15.42578125,19.5703125,18.5703125,15.42578125,19.5703125,14.5703125,15.42578125,19.5703125,14.5703125,15.42578125,19.5703125,14.5703125,15.42578125,19.5703125,14.5703125)
for ($i = 0; $i -lt $widths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i] - 0.8333333333333334
}

# ---------------------------------------------------------------------
# 5. Selection moves from C2 to A2.
# ---------------------------------------------------------------------
$ws.Range("A2").Select()
